$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 174 (date 2026-01-09 / serial 46031) with new totals
$ws.Range("B174").Value = 3343
$ws.Range("C174").Value = 847
$ws.Range("D174").Value = 716
$ws.Range("E174").Value = 1780
$ws.Range("F174").Value = 547

# Append new row 175 (date 2026-01-10 / serial 46032)
$ws.Range("A175").Value = 46032
$ws.Range("A175").NumberFormat = $ws.Range("A174").NumberFormat
$ws.Range("B175").Value = 1282
$ws.Range("C175").Value = 110
$ws.Range("D175").Value = 505
$ws.Range("E175").Value = 667
$ws.Range("F175").Value = 0
